$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71-125 down to 72-126.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Cells.Item(71, 1).Value = 9
$ws.Cells.Item(71, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(71, 3).Value = "Metropolitana"

$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
$ws.Cells.Item(71, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item(71, 5).Value = 13
$ws.Cells.Item(71, 6).Value = 100112026
$ws.Cells.Item(71, 7).Value = "Haba"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 43
$ws.Cells.Item(71, 11).Value = 15000
$ws.Cells.Item(71, 12).Value = 16000
$ws.Cells.Item(71, 13).Value = 15488
$ws.Cells.Item(71, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(71, 16).Value = 620
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
